$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.225.87'
$ws.Range("E2").Value = '  -0.65%  '

$ws.Range("D3").Value = '2.269.93'
$ws.Range("E3").Value = '  -0.98%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.62'
$ws.Range("E5").Value = '  -0.36%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.56'
$ws.Range("E6").Value = '  +1.71%  '

$ws.Range("E7").Value = '  -1.17%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.490'
$ws.Range("E9").Value = '  -1.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.22'
$ws.Range("E10").Value = '  -2.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0789'
$ws.Range("E11").Value = '  -1.78%  '

$ws.Range("E12").Value = '  +0.26%  '

$ws.Range("E13").Value = '  +1.96%  '

$ws.Range("D14").Value = '2.622.08'
$ws.Range("E14").Value = '  -0.94%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.70'
$ws.Range("E15").Value = '  +1.15%  '

$ws.Range("D16").Value = '2.264.96'
$ws.Range("E16").Value = '  -1.56%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.791'
$ws.Range("E17").Value = '  -1.57%  '

$ws.Range("D18").Value = '42.108.69'
$ws.Range("E18").Value = '  -0.69%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.30'
$ws.Range("E19").Value = '  -3.22%  '

$ws.Range("D20").Value = '0.0₃0904'
$ws.Range("E20").Value = '  -1.59%  '

$ws.Range("E21").Value = '  -0.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.72'
$ws.Range("E22").Value = '  -0.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.22'
$ws.Range("E23").Value = '  -2.54%  '

$ws.Range("E24").Value = '  +2.62%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.58'
$ws.Range("E25").Value = '  -1.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.54'
$ws.Range("E27").Value = '  -2.50%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.34'
$ws.Range("E28").Value = '  +3.82%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.59'
$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("E30").Value = '  +1.63%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '162.72'
$ws.Range("E31").Value = '  +0.79%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.26'
$ws.Range("E32").Value = '  -1.53%  '

$ws.Range("E33").Value = '  +0.12%  '

$ws.Range("E34").Value = '  +1.41%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.71'
$ws.Range("E35").Value = '  +2.61%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0735'
$ws.Range("E36").Value = '  -2.90%  '

$ws.Range("E37").Value = '  +0.13%  '

$ws.Range("E38").Value = '  -4.05%  '

$ws.Range("E39").Value = '  -1.94%  '

$ws.Range("E40").Value = '  -1.41%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.10'
$ws.Range("E41").Value = '  -1.88%  '

$ws.Range("E42").Value = '  +2.12%  '

$ws.Range("D43").Value = '1.949.15'
$ws.Range("E43").Value = '  -3.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.04'
$ws.Range("E44").Value = '  -2.46%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0280'
$ws.Range("E45").Value = '  -1.63%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.98'
$ws.Range("E46").Value = '  -1.78%  '

$ws.Range("E47").Value = '  -3.29%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.70'
$ws.Range("E48").Value = '  -0.46%  '

$ws.Range("D49").Value = '2.492.99'
$ws.Range("E49").Value = '  -0.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.35'
$ws.Range("E50").Value = '  -0.98%  '

$ws.Range("E51").Value = '  +0.03%  '
